$d = $word.ActiveDocument

# The document opens with three paragraphs in the body:
#   1) "Hi, my name is Bob Brand and I'm an engineer. ..."
#   2) <empty paragraph>
#   3) "Throughout my career, I have had the opportunity ..."
#
# The target revision removes both paragraphs of biography text in their
# entirety (including their paragraph marks), leaving only the single
# empty paragraph that originally separated them, immediately followed
# by the section properties.

function Remove-ParagraphStartingWith($prefix) {
    # Delete the whole paragraph (including its paragraph mark) whose
    # text starts with $prefix. Using Range.Delete() on the paragraph's
    # own Range (rather than clearing Range.Text) removes the paragraph
    # mark too, so the paragraph disappears instead of becoming blank.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

Remove-ParagraphStartingWith("Hi, my name is Bob Brand") | Out-Null
Remove-ParagraphStartingWith("Throughout my career") | Out-Null
